# Add ability to import existing JSON files
# -------------------------------------------------------------
# This script turns the three "this is the Nth key" placeholder rows
# on the "common" sheet into real looking i18n keys (sign-up /
# our-product-heading / footer), makes the header row bold, widens the
# Russian column a bit, and leaves the "common" sheet as the active /
# selected sheet (instead of "package_settings").

$wb = $excel.ActiveWorkbook

$wsCommon = $wb.Worksheets.Item("common")

# ---- replace the placeholder translation keys with real ones ----
$wsCommon.Range("A2").Value = "sign-up"
$wsCommon.Range("B2").Value = "Sign Up"
$wsCommon.Range("C2").Value = "Регистрация"

$wsCommon.Range("A3").Value = "our-product-heading"
$wsCommon.Range("B3").Value = "This is our great product"
$wsCommon.Range("C3").Value = "Это наш великолепный продукт"

$wsCommon.Range("A4").Value = "footer"
$wsCommon.Range("B4").Value = "Footer"
$wsCommon.Range("C4").Value = "Футер"

# ---- bold the header row (key / en / ru) ----
$wsCommon.Range("A1:C1").Font.Bold = $true

# ---- widen the "ru" column a little ----
$wsCommon.Columns.Item(3).ColumnWidth = 27.87

# ---- make "common" the active sheet, with A9 selected ----
$wsCommon.Select() | Out-Null
$wsCommon.Range("A9").Select() | Out-Null
